# Apply the changes described by the diff:
#  1. Sheet "sets": row 36 (set_id 35, match_id 11, set_number 1) final score
#     changes from 25-20 to 27-21 (home_points / away_points).
#  2. Sheet "rallies": three new rally rows (249-251) are appended, extending
#     the used range from A1:P248 to A1:P251.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update the final set score on the "sets" sheet
# ---------------------------------------------------------------------------
$sets = $wb.Worksheets.Item("sets")
$sets.Cells.Item(36, 4).Value = 27   # D36 home_points: 25 -> 27
$sets.Cells.Item(36, 5).Value = 21   # E36 away_points: 20 -> 21

# ---------------------------------------------------------------------------
# 2) Append new rally rows to the "rallies" sheet
# ---------------------------------------------------------------------------
$rallies = $wb.Worksheets.Item("rallies")

# Column layout:
# A rally_id, B match_id, C set_number, D rally_no, E side, F position,
# G player_number, H action, I result, J who_scored, K score_home,
# L score_away, M raw_text, N position_zone, O pos_fb, P frente_fundo

$newRows = @(
    @{ Row=249; A=258; B=11; C=1; D=46; E="NOS"; G=4; H="LOB";  I="PONTO"; J="NOS"; K=26; L=20; M="1 4 lob";  N="FRENTE"; O="FRENTE"; P="FRENTE" },
    @{ Row=250; A=259; B=11; C=1; D=47; E="NOS"; G=6; H="PIPE"; I="PONTO"; J="NOS"; K=27; L=20; M="1 6 pi";   N="FRENTE"; O="FRENTE"; P="FRENTE" },
    @{ Row=251; A=260; B=11; C=1; D=48; E="NOS"; G=6; H="PIPE"; I="ERRO";  J="ADV"; K=27; L=21; M="1 6 pi e"; N="FRENTE"; O="FRENTE"; P="FRENTE" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $rallies.Cells.Item($row, 1).Value = $r.A
    $rallies.Cells.Item($row, 2).Value = $r.B
    $rallies.Cells.Item($row, 3).Value = $r.C
    $rallies.Cells.Item($row, 4).Value = $r.D
    $rallies.Cells.Item($row, 5).Value = $r.E
    # Column F (position) is an empty string for every data row in this
    # sheet; replicate that by copying an existing empty-string cell into
    # the new row instead of assigning "" directly (which would clear the
    # cell entirely rather than leaving an empty text value).
    $rallies.Cells.Item(2, 6).Copy($rallies.Cells.Item($row, 6))
    $rallies.Cells.Item($row, 7).Value = $r.G
    $rallies.Cells.Item($row, 8).Value = $r.H
    $rallies.Cells.Item($row, 9).Value = $r.I
    $rallies.Cells.Item($row, 10).Value = $r.J
    $rallies.Cells.Item($row, 11).Value = $r.K
    $rallies.Cells.Item($row, 12).Value = $r.L
    $rallies.Cells.Item($row, 13).Value = $r.M
    $rallies.Cells.Item($row, 14).Value = $r.N
    $rallies.Cells.Item($row, 15).Value = $r.O
    $rallies.Cells.Item($row, 16).Value = $r.P
}
